# Applies the "Fake user data added to templates.xslx" commit:
#  - selects F4:G4 on the "Template" sheet (was A7)
#  - selects D15:E15 on the "Points of interest" sheet (was D20)
#  - adds a new "Fake Points - Florida" sheet at the end, populated with
#    fake per-user lat/long data, which becomes the active sheet/tab

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Template sheet: selection moves from A7 to F4:G4
# ---------------------------------------------------------------------
$wsTemplate = $wb.Worksheets.Item("Template")
$wsTemplate.Range("F4:G4").Select() | Out-Null

# ---------------------------------------------------------------------
# 2. Points of interest sheet: selection moves from D20 to D15:E15
# ---------------------------------------------------------------------
$wsPoints = $wb.Worksheets.Item("Points of interest")
$wsPoints.Range("D15:E15").Select() | Out-Null

# ---------------------------------------------------------------------
# 3. New "Fake Points - Florida" sheet, appended after the last sheet
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws3 = $wb.Worksheets.Add([System.Type]::Missing, $lastSheet)
$ws3.Name = "Fake Points - Florida"

# Header row. Columns B and C intentionally share the same label
# ("Lead User Latitude") exactly as in the source data.
$ws3.Range("A1").Value = "instance"
$ws3.Range("B1").Value = "Lead User Latitude"
$ws3.Range("C1").Value = "Lead User Latitude"
$ws3.Range("D1").Value = "User 2 Latitude"
$ws3.Range("F1").Value = "User 3 Latitude"
$ws3.Range("E1").Value = "User 2 Longitude"
$ws3.Range("G1").Value = "User 3 Longitude"

# Data rows 2-15: instance id + lead/user2/user3 lat-long samples.
$data = @(
  @(1, 29.9107507079447, -81.313592409669596, 29.911420290532501, -81.312390780030796, 29.9120898686188, -81.311704134522699),
  @(2, 29.835431168878301, -81.320716356812696, 29.842578546534899, -81.322089647828804, 29.847343180827899, -81.323462938844003),
  @(3, 29.549714853149499, -81.285010790406304, 29.577188399379398, -81.286384081421403, 29.610624366493902, -81.291877245484002),
  @(4, 29.345221060567901, -81.149054979859002, 29.355994157883099, -81.162787890015295, 29.365569288019799, -81.172400927125196),
  @(5, 29.152317645156401, -81.0666575189219, 29.195484563183001, -81.088630175171104, 29.219458336307099, -81.105109667358704),
  @(6, 28.830394634201799, -80.907355761109301, 28.859264379655301, -80.9183420892344, 29.113931815349002, -81.039191698609201),
  @(7, 28.830394634201799, -80.907355761109301, 28.859264379655301, -80.9183420892344, 28.883316378915101, -80.9293284173595),
  @(8, 28.459190425189401, -80.830451464232993, 28.386725332899601, -80.802985643921204, 28.4495312793991, -80.813971972046303),
  @(9, 28.106061397387599, -80.709601854858093, 28.183559669324701, -80.715095018920707, 28.193243008970502, -80.709601854858093),
  @(10, 27.4304486507646, -80.407477831420593, 27.5279169782251, -80.456916307983406, 27.557140649474899, -80.489875292358604),
  @(11, 27.163799994017499, -80.305854296263703, 27.190676883070701, -80.327826952513902, 27.215104799611002, -80.344306444701502),
  @(12, 27.0562277718441, -80.121833300169797, 27.063565524779001, -80.116340136108093, 27.0684570933517, -80.1135935540768),
  @(13, 26.999955725568999, -80.091620897826701, 27.034211633332699, -80.105353807983093, 27.0684570933517, -80.1135935540768),
  @(14, 26.952224746152101, -80.090247606811502, 26.984047648061001, -80.094367479857993, 26.977928557744999, -80.091620897826701)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $rowNum = $i + 2
    $rowVals = $data[$i]
    for ($j = 0; $j -lt $rowVals.Length; $j++) {
        $ws3.Cells.Item($rowNum, $j + 1).Value = $rowVals[$j]
    }
}

# Annotation notes in column I on a couple of rows.
$ws3.Range("I7").Value = "*** User Three Encount Problem"
$ws3.Range("I12").Value = "*** Group decides to go to Juptier Island"

# Trailing rows 16-20 just continue the instance counter in column A.
$ws3.Range("A16").Value = 15
$ws3.Range("A17").Value = 16
$ws3.Range("A18").Value = 17
$ws3.Range("A19").Value = 18
$ws3.Range("A20").Value = 19

# Column widths (best-fit-like sizes from the source workbook).
$ws3.Columns.Item(1).ColumnWidth = 8.42578125
$ws3.Columns.Item(2).ColumnWidth = 17.5703125
$ws3.Columns.Item(3).ColumnWidth = 17.5703125
$ws3.Columns.Item(4).ColumnWidth = 14.28515625
$ws3.Columns.Item(5).ColumnWidth = 15.85546875
$ws3.Columns.Item(6).ColumnWidth = 14.28515625
$ws3.Columns.Item(7).ColumnWidth = 15.85546875

# Final selection/active cell on the new sheet.
$ws3.Range("B16").Select() | Out-Null
